$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (e.g. "241.68") must be
# forced to Text format first, otherwise Excel auto-converts them to numeric
# values instead of keeping the original text representation.
$numericLookingCells = @(
    "D6", "D8", "D9", "D10", "D11", "D12", "D13", "D15", "D17", "D18", "D19", "D20", "D23", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D33", "D34", "D35", "D36", "D37", "D38", "D40", "D41", "D42", "D43", "D45", "D47", "D48", "D49", "D51"
)
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '29.503.50'
$ws.Range("E2").Value = '  +1.07%  '
$ws.Range("D3").Value = '1.878.75'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("E5").Value = '  +2.56%  '
$ws.Range("D6").Value = '241.68'
$ws.Range("E6").Value = '  +1.80%  '
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = '0.07927'
$ws.Range("E8").Value = '  +1.11%  '
$ws.Range("D9").Value = '0.3108'
$ws.Range("E9").Value = '  +3.17%  '
$ws.Range("D10").Value = '25.22'
$ws.Range("E10").Value = '  +6.34%  '
$ws.Range("D11").Value = '0.08278'
$ws.Range("E11").Value = '  +2.00%  '
$ws.Range("D12").Value = '0.7291'
$ws.Range("E12").Value = '  +3.55%  '
$ws.Range("D13").Value = '5.282'
$ws.Range("E13").Value = '  +2.07%  '
$ws.Range("D14").Value = '1.854.03'
$ws.Range("E14").Value = '  +0.15%  '
$ws.Range("D15").Value = '91.15'
$ws.Range("E15").Value = '  +1.98%  '
$ws.Range("D16").Value = '29.497.47'
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("D17").Value = '5.932'
$ws.Range("E17").Value = '  +2.41%  '
$ws.Range("D18").Value = '246.32'
$ws.Range("D19").Value = '0.000007882'
$ws.Range("E19").Value = '  +0.90%  '
$ws.Range("D20").Value = '13.34'
$ws.Range("E20").Value = '  +1.13%  '
$ws.Range("D21").Value = '2.125.35'
$ws.Range("E21").Value = '  +1.15%  '
$ws.Range("E22").Value = '  +0.06%  '
$ws.Range("D23").Value = '7.955'
$ws.Range("E23").Value = '  +6.06%  '
$ws.Range("E24").Value = '  +0.08%  '
$ws.Range("D25").Value = '0.1625'
$ws.Range("E25").Value = '  +15.05%  '
$ws.Range("D26").Value = '163.59'
$ws.Range("E26").Value = '  +0.66%  '
$ws.Range("D27").Value = '9.060'
$ws.Range("E27").Value = '  +2.23%  '
$ws.Range("E28").Value = '  +1.86%  '
$ws.Range("D29").Value = '1.360'
$ws.Range("E29").Value = '  -2.94%  '
$ws.Range("D30").Value = '1.493'
$ws.Range("E30").Value = '  +1.22%  '
$ws.Range("D31").Value = '4.386'
$ws.Range("E31").Value = '  +1.73%  '
$ws.Range("D32").Value = '4.111'
$ws.Range("E32").Value = '  +2.61%  '
$ws.Range("D33").Value = '0.05269'
$ws.Range("E33").Value = '  +2.43%  '
$ws.Range("D34").Value = '1.960'
$ws.Range("E34").Value = '  +2.23%  '
$ws.Range("D35").Value = '1.199'
$ws.Range("E35").Value = '  +2.83%  '
$ws.Range("D36").Value = '0.7271'
$ws.Range("E36").Value = '  +2.65%  '
$ws.Range("D37").Value = '2.680'
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("D38").Value = '0.01868'
$ws.Range("E38").Value = '  +1.34%  '
$ws.Range("D39").Value = '1.229.24'
$ws.Range("E39").Value = '  +6.85%  '
$ws.Range("D40").Value = '2.719'
$ws.Range("E40").Value = '  +0.61%  '
$ws.Range("D41").Value = '0.9122'
$ws.Range("E41").Value = '  -1.06%  '
$ws.Range("B42").Value = 'Aave'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D42").Value = '73.74'
$ws.Range("E42").Value = '  +5.32%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '6.148'
$ws.Range("E43").Value = '  +3.20%  '
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").Value = '102.07'
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").Value = '2.020.14'
$ws.Range("E46").Value = '  +1.12%  '
$ws.Range("D47").Value = '0.5291'
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("D48").Value = '1.803'
$ws.Range("E48").Value = '  +4.06%  '
$ws.Range("D49").Value = '2.931'
$ws.Range("E49").Value = '  +10.54%  '
$ws.Range("E50").Value = '  +1.73%  '
$ws.Range("D51").Value = '9.325'
$ws.Range("E51").Value = '  +2.04%  '

# Restore the default (unstyled) cell style on the cells we temporarily
# switched to Text number format, so the resulting style matches the original.
foreach ($addr in $numericLookingCells) {
    $ws.Range($addr).Style = "Normal"
}
